$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose textual values must be protected from Excel's automatic
# type-inference (dates like "2023-09-16" or digit-only strings like "3"
# must remain plain text, exactly like the rest of this sheet stores them).
$textProtectCols = @("Y", "AA")

function Set-TextCell($col, $row, $value) {
    $addr = "$col$row"
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

function Set-NumCell($col, $row, $value) {
    $addr = "$col$row"
    $ws.Range($addr).Value = $value
}

function Set-BoolCell($col, $row, $value) {
    $addr = "$col$row"
    $ws.Range($addr).Value = $value
}

$rows = @(
    @{
        Row = 54
        Num = @{ A = 112128672; B = 90689; E = 5966; Q = 623279.5584157004; R = 6951830.900261804; S = 100 }
        Txt = @{
            C = "Ovaliderad"; D = "NT"; F = "Motaggsvamp"; G = "Sarcodon squamosus";
            H = "(Schaeff.) Quél."; P = "Mjällådalen, Mpd"; T = "Västernorrland";
            U = "Timrå"; V = "Medelpad"; W = "Ljustorp";
            Y = "2023-09-16"; Z = "00:00"; AA = "2023-09-16"; AB = "00:00";
            AC = "stgen upp mot vägen";
            AW = "Elisabeth Nilsson"; AX = "Elisabeth Nilsson"
        }
        Bool = @{ AD = $false; AE = $false; AG = $false }
    },
    @{
        Row = 55
        Num = @{ A = 112128712; B = 88914; E = 2051; Q = 623279.5584157004; R = 6951830.900261804; S = 100 }
        Txt = @{
            C = "Ovaliderad"; D = "VU"; F = "Rotfingersvamp"; G = "Ramaria boreimaxima";
            H = "Kytöv. & M.Toivonen"; I = "3"; J = "fruktkroppar";
            P = "Mjällådalen, Mpd"; T = "Västernorrland";
            U = "Timrå"; V = "Medelpad"; W = "Ljustorp";
            Y = "2023-09-16"; Z = "00:00"; AA = "2023-09-16"; AB = "00:00";
            AC = "stigen upp mot vägen";
            AW = "Elisabeth Nilsson"; AX = "Elisabeth Nilsson"
        }
        Bool = @{ AD = $false; AE = $false; AG = $false }
    },
    @{
        Row = 56
        Num = @{ A = 112128664; B = 90710; E = 5449; Q = 623279.5584157004; R = 6951830.900261804; S = 100 }
        Txt = @{
            C = "Ovaliderad"; D = "NT"; F = "Svart taggsvamp"; G = "Phellodon niger";
            H = "(Fr.:Fr.) P.Karst."; P = "Mjällådalen, Mpd"; T = "Västernorrland";
            U = "Timrå"; V = "Medelpad"; W = "Ljustorp";
            Y = "2023-09-16"; Z = "00:00"; AA = "2023-09-16"; AB = "00:00";
            AC = "stigen upp mot vägen";
            AW = "Elisabeth Nilsson"; AX = "Elisabeth Nilsson"
        }
        Bool = @{ AD = $false; AE = $false; AG = $false }
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    foreach ($col in $r.Num.Keys) {
        Set-NumCell $col $rowNum $r.Num[$col]
    }

    foreach ($col in $r.Txt.Keys) {
        if ($textProtectCols -contains $col -or $col -eq "I") {
            Set-TextCell $col $rowNum $r.Txt[$col]
        } else {
            Set-NumCell $col $rowNum $r.Txt[$col]
        }
    }

    foreach ($col in $r.Bool.Keys) {
        Set-BoolCell $col $rowNum $r.Bool[$col]
    }
}
